# Merge and unmerge cells using openpyxl
# -----------------------------------------------------------------
# This script reproduces, via Excel COM automation, the structural
# edit that the original author made with openpyxl:
#   * Added a "Sort" worksheet containing two columns of unsorted
#     random integers (data to demonstrate a later sort exercise).
#   * Added a "Merge" worksheet containing two merged cell ranges
#     (and the corresponding "unmerge" companion data).
#   * The newly-added "Merge" sheet becomes the active / selected
#     worksheet tab.
# -----------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Add the "Sort" worksheet (placed after the last existing sheet)
# ---------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sortSheet = $wb.Worksheets.Add($null, $lastSheet)
$sortSheet.Name = "Sort"

$sortSheet.Columns.Item(1).ColumnWidth = 10.45

$sortData = @(
    @(40, 26),
    @(93, 45),
    @(23, 54),
    @(80, 43),
    @(21, 12),
    @(63, 29),
    @(34, 15),
    @(80, 68),
    @(20, 41)
)

$r = 3
foreach ($row in $sortData) {
    $sortSheet.Cells.Item($r, 1).Value = $row[0]
    $sortSheet.Cells.Item($r, 2).Value = $row[1]
    $r++
}

# Select the whole sheet (matches the "select all before sort" view
# state left behind by the authoring tool).
$sortSheet.Cells.Select() | Out-Null

# ---------------------------------------------------------------
# 2) Add the "Merge" worksheet (placed after "Sort")
# ---------------------------------------------------------------
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$mergeSheet = $wb.Worksheets.Add($null, $lastSheet2)
$mergeSheet.Name = "Merge"

# Merge a header-style band and a larger block lower on the sheet.
$mergeSheet.Range("A3:E3").Merge() | Out-Null
$mergeSheet.Range("F9:J15").Merge() | Out-Null

# Touch the intervening rows so the sheet's used range/dimension
# spans down through row 15 (matches the target row extent).
foreach ($rowIdx in @(3, 9, 10, 11, 12, 13, 14, 15)) {
    $mergeSheet.Rows.Item($rowIdx).OutlineLevel = 0
}

$mergeSheet.Range("A3:E3").Select() | Out-Null

Write-Output "Added Sort and Merge worksheets"
